$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.309.16'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.790.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +1.07%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.62'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0688'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0944'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.051.59'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.798.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.632'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.373.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0792'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '243.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '166.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.99%  '
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.115'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0524'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.80'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E35').Value = '  -2.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.397.52'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.670'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.13'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.89%  '
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.933'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0525'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.58%  '
$ws.Range('E46').Value = '  +2.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.01'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.951.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('E51').Value = '  -1.85%  '
